$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 116
$ws.Range("H116").Value = 7389.3076
$ws.Range("I116").Value = 3935.6428
$ws.Range("J116").Value = 11418.583
$ws.Range("K116").Value = 3935.6428
$ws.Range("L116").Value = 11418.583
$ws.Range("M116").Value = -493.6428000000001
$ws.Range("N116").Value = -18302.583

# Row 132
$ws.Range("H132").Value = 2518.9697
$ws.Range("I132").Value = 2302.6924
$ws.Range("K132").Value = 6908.0772
$ws.Range("M132").Value = -4378.0772

# Row 133
$ws.Range("H133").Value = 38737.527
$ws.Range("J133").Value = 38737.527
$ws.Range("L133").Value = 38737.527
$ws.Range("N133").Value = -48857.527

# Row 136
$ws.Range("H136").Value = 48333.332
$ws.Range("J136").Value = 48333.332
$ws.Range("L136").Value = 48333.332
$ws.Range("N136").Value = -58533.332

# Row 138
$ws.Range("H138").Value = 4203.804
$ws.Range("I138").Value = 2959.4443
$ws.Range("J138").Value = 4506.4863
$ws.Range("K138").Value = 8878.332900000001
$ws.Range("L138").Value = 13519.4589
$ws.Range("M138").Value = -3738.332900000001
$ws.Range("N138").Value = -23799.4589

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3688.1667
$ws.Range("I32").Value = 3688.1667
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3688.1667
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -3401.1667
$ws.Range("N32").ClearContents()

# Row 61
$ws.Range("H61").Value = 2400.611
$ws.Range("I61").Value = 2394.7646
$ws.Range("K61").Value = 2394.7646
$ws.Range("M61").Value = -2182.7646

# Row 74
$ws.Range("H74").Value = 1385.4634
$ws.Range("I74").Value = 1312.625
$ws.Range("J74").Value = 1644.4445
$ws.Range("K74").Value = 1312.625
$ws.Range("L74").Value = 1644.4445
$ws.Range("M74").Value = -438.625
$ws.Range("N74").Value = -3392.4445

# Row 77
$ws.Range("H77").Value = 1385.4634
$ws.Range("I77").Value = 1312.625
$ws.Range("J77").Value = 1644.4445
$ws.Range("K77").Value = 6563.125
$ws.Range("L77").Value = 8222.2225
$ws.Range("M77").Value = -2195.125
$ws.Range("N77").Value = -16958.2225

# Row 122
$ws.Range("H122").Value = 4083.037
$ws.Range("I122").Value = 2260.923
$ws.Range("J122").Value = 5775
$ws.Range("K122").Value = 6782.768999999999
$ws.Range("L122").Value = 17325
$ws.Range("M122").Value = -4332.768999999999
$ws.Range("N122").Value = -22225

# Row 136
$ws.Range("H136").Value = 2400.611
$ws.Range("I136").Value = 2394.7646
$ws.Range("K136").Value = 7184.293799999999
$ws.Range("M136").Value = -4634.293799999999

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 557979.25
$ws.Range("I107").Value = 2348.0833
$ws.Range("J107").Value = 1669241.6
$ws.Range("K107").Value = 2348.0833
$ws.Range("L107").Value = 1669241.6
$ws.Range("M107").Value = -428.0832999999998
$ws.Range("N107").Value = -1673081.6

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 129975.625
$ws.Range("I31").Value = 1453
$ws.Range("J31").Value = 258498.25
$ws.Range("K31").Value = 1453
$ws.Range("L31").Value = 258498.25
$ws.Range("M31").Value = -1158
$ws.Range("N31").Value = -259088.25

# Row 32
$ws.Range("H32").Value = 337333.34
$ws.Range("I32").Value = 337333.34
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 337333.34
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -337017.34
$ws.Range("N32").ClearContents()

# Row 34
$ws.Range("H34").Value = 129975.625
$ws.Range("I34").Value = 1453
$ws.Range("J34").Value = 258498.25
$ws.Range("K34").Value = 1453
$ws.Range("L34").Value = 258498.25
$ws.Range("M34").Value = -1251
$ws.Range("N34").Value = -258902.25

# Row 132
$ws.Range("H132").Value = 2286.261
$ws.Range("I132").Value = 1844.7333
$ws.Range("J132").Value = 3114.125
$ws.Range("K132").Value = 5534.199900000001
$ws.Range("L132").Value = 9342.375
$ws.Range("M132").Value = -3004.199900000001
$ws.Range("N132").Value = -14402.375

$ws = $wb.Worksheets.Item("CUL")
# Row 81
$ws.Range("H81").Value = 208416.28
$ws.Range("J81").Value = 208416.28
$ws.Range("L81").Value = 625248.84
$ws.Range("N81").Value = -627494.84

# Row 84
$ws.Range("H84").Value = 208416.28
$ws.Range("J84").Value = 208416.28
$ws.Range("L84").Value = 1875746.52
$ws.Range("N84").Value = -1886978.52

# Row 129
$ws.Range("H129").Value = 7048.294
$ws.Range("I129").Value = 448
$ws.Range("J129").Value = 19148.834
$ws.Range("K129").Value = 1344
$ws.Range("L129").Value = 57446.50199999999
$ws.Range("M129").Value = 3656
$ws.Range("N129").Value = -67446.50199999999

# Row 133
$ws.Range("H133").Value = 17732.521
$ws.Range("J133").Value = 21235.234
$ws.Range("L133").Value = 63705.702
$ws.Range("N133").Value = -73825.702

# Row 139
$ws.Range("H139").Value = 7747.6
$ws.Range("I139").Value = 5107
$ws.Range("J139").Value = 8153.846
$ws.Range("K139").Value = 15321
$ws.Range("L139").Value = 24461.538
$ws.Range("M139").Value = -10181
$ws.Range("N139").Value = -34741.538

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 3684.1904
$ws.Range("I102").Value = 1861.0714
$ws.Range("K102").Value = 1861.0714
$ws.Range("M102").Value = -239.0714

# Row 126
$ws.Range("H126").Value = 3004.6924
$ws.Range("I126").Value = 1013.8
$ws.Range("K126").Value = 3041.4
$ws.Range("M126").Value = -571.3999999999996

# Row 141
$ws.Range("H141").Value = 48000
$ws.Range("J141").Value = 48000
$ws.Range("L141").Value = 48000
$ws.Range("N141").Value = -58360

$ws = $wb.Worksheets.Item("LTW")
# Row 32
$ws.Range("H32").Value = 2000
$ws.Range("I32").Value = 2000
$ws.Range("K32").Value = 2000
$ws.Range("M32").Value = -1683

# Row 136
$ws.Range("H136").Value = 291884.5
$ws.Range("I136").Value = 504340.06
$ws.Range("J136").Value = 8610.467000000001
$ws.Range("K136").Value = 1513020.18
$ws.Range("L136").Value = 25831.401
$ws.Range("M136").Value = -1510470.18
$ws.Range("N136").Value = -30931.401

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 1378.3077
$ws.Range("I107").Value = 1545.0952
$ws.Range("J107").Value = 677.8
$ws.Range("K107").Value = 4635.2856
$ws.Range("L107").Value = 2033.4
$ws.Range("M107").Value = -2715.2856
$ws.Range("N107").Value = -5873.4

# Row 135
$ws.Range("H135").Value = 60800
$ws.Range("J135").Value = 60800
$ws.Range("L135").Value = 60800
$ws.Range("N135").Value = -70940

# Row 141
$ws.Range("H141").Value = 53333.332
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 53333.332
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 53333.332
$ws.Range("N141").Value = -63693.332
$ws.Range("M141").ClearContents()
